# Scheduled runner: refresh market-price-derived profit columns
# (currentAveragePrice[NQ/HQ], LevePrice[NQ/HQ], LeveProfit[NQ/HQ]) for
# the rows whose backing Universalis price lookups changed.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 558.29785
$ws.Range("I33").Value = 374.14633
$ws.Range("J33").Value = 1816.6666
$ws.Range("K33").Value = 374.14633
$ws.Range("L33").Value = 1816.6666
$ws.Range("M33").Value = -145.14633
$ws.Range("N33").Value = -2274.6666
$ws.Range("H58").Value = 1180.4117
$ws.Range("J58").Value = 1285.4286
$ws.Range("L58").Value = 3856.2858
$ws.Range("N58").Value = -4156.2858
$ws.Range("H137").Value = 3005.8865
$ws.Range("I137").Value = 2186.5625
$ws.Range("J137").Value = 5190.75
$ws.Range("K137").Value = 6559.6875
$ws.Range("L137").Value = 15572.25
$ws.Range("M137").Value = -4009.6875
$ws.Range("N137").Value = -20672.25
$ws.Range("H138").Value = 2250.7344
$ws.Range("I138").Value = 2927.5715
$ws.Range("J138").Value = 2061.22
$ws.Range("K138").Value = 8782.7145
$ws.Range("L138").Value = 6183.66
$ws.Range("M138").Value = -3642.7145
$ws.Range("N138").Value = -16463.66

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2009.1818
$ws.Range("I2").Value = 2009.1818
$ws.Range("K2").Value = 2009.1818
$ws.Range("M2").Value = -1896.1818
$ws.Range("H32").Value = 436206.06
$ws.Range("I32").Value = 492547.72
$ws.Range("K32").Value = 492547.72
$ws.Range("M32").Value = -492260.72
$ws.Range("H61").Value = 2049.625
$ws.Range("I61").Value = 1228.25
$ws.Range("J61").Value = 3966.1667
$ws.Range("K61").Value = 1228.25
$ws.Range("L61").Value = 3966.1667
$ws.Range("M61").Value = -1016.25
$ws.Range("N61").Value = -4390.1667
$ws.Range("H74").Value = 1824.2727
$ws.Range("I74").Value = 1348.6072
$ws.Range("J74").Value = 2656.6875
$ws.Range("K74").Value = 1348.6072
$ws.Range("L74").Value = 2656.6875
$ws.Range("M74").Value = -474.6071999999999
$ws.Range("N74").Value = -4404.6875
$ws.Range("H77").Value = 1824.2727
$ws.Range("I77").Value = 1348.6072
$ws.Range("J77").Value = 2656.6875
$ws.Range("K77").Value = 6743.036
$ws.Range("L77").Value = 13283.4375
$ws.Range("M77").Value = -2375.036
$ws.Range("N77").Value = -22019.4375
$ws.Range("H116").Value = 2009.1818
$ws.Range("I116").Value = 2009.1818
$ws.Range("K116").Value = 2009.1818
$ws.Range("M116").Value = 284.8181999999999
$ws.Range("H136").Value = 2049.625
$ws.Range("I136").Value = 1228.25
$ws.Range("J136").Value = 3966.1667
$ws.Range("K136").Value = 3684.75
$ws.Range("L136").Value = 11898.5001
$ws.Range("M136").Value = -1134.75
$ws.Range("N136").Value = -16998.5001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2009.1818
$ws.Range("I3").Value = 2009.1818
$ws.Range("K3").Value = 2009.1818
$ws.Range("M3").Value = -1895.1818
$ws.Range("H52").Value = 44500
$ws.Range("J52").Value = 44500
$ws.Range("L52").Value = 44500
$ws.Range("N52").Value = -45026
$ws.Range("H55").Value = 42750
$ws.Range("J55").Value = 42750
$ws.Range("L55").Value = 42750
$ws.Range("N55").Value = -43296
$ws.Range("H121").Value = 44500
$ws.Range("J121").Value = 44500
$ws.Range("L121").Value = 44500
$ws.Range("N121").Value = -47994
$ws.Range("H134").Value = 2679.8928
$ws.Range("I134").Value = 2581.15
$ws.Range("J134").Value = 2926.75
$ws.Range("K134").Value = 7743.450000000001
$ws.Range("L134").Value = 8780.25
$ws.Range("M134").Value = -5208.450000000001
$ws.Range("N134").Value = -13850.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5724.569
$ws.Range("I31").Value = 1486.2593
$ws.Range("J31").Value = 10492.667
$ws.Range("K31").Value = 1486.2593
$ws.Range("L31").Value = 10492.667
$ws.Range("M31").Value = -1191.2593
$ws.Range("N31").Value = -11082.667
$ws.Range("H34").Value = 5724.569
$ws.Range("I34").Value = 1486.2593
$ws.Range("J34").Value = 10492.667
$ws.Range("K34").Value = 1486.2593
$ws.Range("L34").Value = 10492.667
$ws.Range("M34").Value = -1284.2593
$ws.Range("N34").Value = -10896.667
$ws.Range("H58").Value = 1052.7297
$ws.Range("I58").Value = 747.6087
$ws.Range("K58").Value = 747.6087
$ws.Range("M58").Value = -544.6087
$ws.Range("H107").Value = 3906888.5
$ws.Range("I107").Value = 6250561.5
$ws.Range("J107").Value = 766.6667
$ws.Range("K107").Value = 6250561.5
$ws.Range("L107").Value = 766.6667
$ws.Range("M107").Value = -6248641.5
$ws.Range("N107").Value = -4606.6667
$ws.Range("H132").Value = 9260986
$ws.Range("I132").Value = 1121.9166
$ws.Range("K132").Value = 3365.7498
$ws.Range("M132").Value = -835.7498000000001
$ws.Range("H134").Value = 2209.6428
$ws.Range("I134").Value = 1394.5652
$ws.Range("K134").Value = 4183.6956
$ws.Range("M134").Value = -1648.6956
$ws.Range("H136").Value = 1052.7297
$ws.Range("I136").Value = 747.6087
$ws.Range("K136").Value = 2242.8261
$ws.Range("M136").Value = 307.1738999999998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 520.2857
$ws.Range("I5").Value = 440.33334
$ws.Range("K5").Value = 1321.00002
$ws.Range("M5").Value = -1209.00002
$ws.Range("H88").Value = 10845.25
$ws.Range("J88").Value = 10845.25
$ws.Range("L88").Value = 32535.75
$ws.Range("N88").Value = -33391.75
$ws.Range("H91").Value = 10845.25
$ws.Range("J91").Value = 10845.25
$ws.Range("L91").Value = 32535.75
$ws.Range("N91").Value = -35499.75
$ws.Range("H122").Value = 6219.3335
$ws.Range("I122").Value = 376.5
$ws.Range("J122").Value = 13522.875
$ws.Range("K122").Value = 3388.5
$ws.Range("L122").Value = 121705.875
$ws.Range("M122").Value = -938.5
$ws.Range("N122").Value = -126605.875
$ws.Range("H131").Value = 983.8182
$ws.Range("J131").Value = 1110.7778
$ws.Range("L131").Value = 3332.3334
$ws.Range("N131").Value = -13412.3334
$ws.Range("H135").Value = 520.2857
$ws.Range("I135").Value = 440.33334
$ws.Range("K135").Value = 3963.00006
$ws.Range("M135").Value = -1428.00006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2121.611
$ws.Range("J126").Value = 2430.5715
$ws.Range("L126").Value = 7291.7145
$ws.Range("N126").Value = -12231.7145

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1025
$ws.Range("I46").Value = 957.1429000000001
$ws.Range("J46").Value = 1500
$ws.Range("K46").Value = 957.1429000000001
$ws.Range("L46").Value = 1500
$ws.Range("M46").Value = -769.1429000000001
$ws.Range("N46").Value = -1876
$ws.Range("H55").Value = 798.73334
$ws.Range("I55").Value = 433.66666
$ws.Range("J55").Value = 890
$ws.Range("K55").Value = 433.66666
$ws.Range("L55").Value = 890
$ws.Range("M55").Value = -260.66666
$ws.Range("N55").Value = -1236
$ws.Range("H132").Value = 4400.7334
$ws.Range("I132").Value = 3601.9
$ws.Range("K132").Value = 10805.7
$ws.Range("M132").Value = -8275.700000000001
$ws.Range("H136").Value = 3704976.5
$ws.Range("I136").Value = 918.3871
$ws.Range("J136").Value = 11906819
$ws.Range("K136").Value = 2755.1613
$ws.Range("L136").Value = 35720457
$ws.Range("M136").Value = -205.1613000000002
$ws.Range("N136").Value = -35725557

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 7971
$ws.Range("J74").Value = 7565.2
$ws.Range("L74").Value = 7565.2
$ws.Range("N74").Value = -9437.200000000001
$ws.Range("H77").Value = 7971
$ws.Range("J77").Value = 7565.2
$ws.Range("L77").Value = 22695.6
$ws.Range("N77").Value = -32055.6
$ws.Range("H132").Value = 3878288
$ws.Range("I132").Value = 2758.0588
$ws.Range("K132").Value = 8274.1764
$ws.Range("M132").Value = -5744.1764
$ws.Range("H136").Value = 2691.7932
$ws.Range("I136").Value = 2335.9167
$ws.Range("J136").Value = 4400
$ws.Range("K136").Value = 7007.750100000001
$ws.Range("L136").Value = 13200
$ws.Range("M136").Value = -4457.750100000001
$ws.Range("N136").Value = -18300

Write-Output "Applied 200 cell updates"